# "warning + pref tab update" — add two new header columns to the
# InformationSheet table:
#   * a new first column "Recommendations For Use" (before "For Best Results")
#   * a new column "Ingredients" (between "Our promises to you" and
#     "Used By & Best Before Date")
# The existing data columns keep their row-2 values; the two brand-new
# columns stay blank in row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new first column (A). Everything that was A:E shifts to B:F.
$ws.Columns("A").Insert()

# 2) Insert a new column at F (currently holding the old "Used By & Best
#    Before Date" column after the first shift). It shifts to G.
$ws.Columns("F").Insert()

# 3) Fill in the two new header cells.
$ws.Range("A1").Value = "Recommendations For Use"
$ws.Range("F1").Value = "Ingredients"

# 4) New columns come in with default formatting; match the bold / bordered /
#    centered header style used by the rest of row 1 (style copies only the
#    formatting, not the value, of the neighboring header cell).
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("G1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Row 2 under the two new columns (A2, F2) is intentionally left blank,
# matching the target layout.
